$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the old row 147, shifting the existing
# rows 147:181 down to 149:183 (dimension grows from A1:T181 to A1:T183).
$ws.Rows("147:148").Insert()

# Fill the first new row (147) - Chirimoya, Primera, week of 2022-08-12 (serial 44785)
$ws.Cells.Item(147, 1).Value = 3
$ws.Cells.Item(147, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(147, 3).Value = "Coquimbo"
$ws.Cells.Item(147, 4).Value = 44785
$ws.Cells.Item(147, 5).Value = 5
$ws.Cells.Item(147, 6).Value = "Fruta"
$ws.Cells.Item(147, 7).Value = 100107
$ws.Cells.Item(147, 8).Value = "Otros"
$ws.Cells.Item(147, 9).Value = 100107002
$ws.Cells.Item(147, 10).Value = "Chirimoya"
$ws.Cells.Item(147, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(147, 12).Value = "Primera"
$ws.Cells.Item(147, 13).Value = 45
$ws.Cells.Item(147, 14).Value = 30000
$ws.Cells.Item(147, 15).Value = 30000
$ws.Cells.Item(147, 16).Value = 30000
$ws.Cells.Item(147, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(147, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(147, 19).Value = 3000
$ws.Cells.Item(147, 20).Value = 10

# Fill the second new row (148) - Chirimoya, Segunda, same week
$ws.Cells.Item(148, 1).Value = 3
$ws.Cells.Item(148, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(148, 3).Value = "Coquimbo"
$ws.Cells.Item(148, 4).Value = 44785
$ws.Cells.Item(148, 5).Value = 5
$ws.Cells.Item(148, 6).Value = "Fruta"
$ws.Cells.Item(148, 7).Value = 100107
$ws.Cells.Item(148, 8).Value = "Otros"
$ws.Cells.Item(148, 9).Value = 100107002
$ws.Cells.Item(148, 10).Value = "Chirimoya"
$ws.Cells.Item(148, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(148, 12).Value = "Segunda"
$ws.Cells.Item(148, 13).Value = 47
$ws.Cells.Item(148, 14).Value = 27000
$ws.Cells.Item(148, 15).Value = 27000
$ws.Cells.Item(148, 16).Value = 27000
$ws.Cells.Item(148, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(148, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(148, 19).Value = 2700
$ws.Cells.Item(148, 20).Value = 10

# Ensure date cells keep the date number format used throughout column D
$ws.Cells.Item(147, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(148, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
